$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '98.474.00'
$ws.Range("E2").Value = '  -0.37%  '

$ws.Range("D3").Value = '3.338.64'
$ws.Range("E3").Value = '  -1.61%  '

$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.15%  '

$ws.Range("D5").Value = '261.75'
$ws.Range("E5").Value = '  +0.32%  '

$ws.Range("D6").Value = '648.73'
$ws.Range("E6").Value = '  +1.92%  '

$ws.Range("E7").Value = '  +8.57%  '

$ws.Range("D8").Value = '0.452'
$ws.Range("E8").Value = '  +13.73%  '

$ws.Range("E9").Value = '  +19.96%  '

$ws.Range("D10").Value = '0.999'
$ws.Range("E10").Value = '  -0.09%  '

$ws.Range("D11").Value = '3.332.43'
$ws.Range("E11").Value = '  -1.65%  '

$ws.Range("B12").Value = 'TRON'
$ws.Range("C12").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D12").Value = '0.207'
$ws.Range("E12").Value = '  +3.34%  '

$ws.Range("B13").Value = 'Avalanche'
$ws.Range("C13").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D13").Value = '43.76'
$ws.Range("E13").Value = '  +20.29%  '

$ws.Range("E14").Value = '  +7.45%  '

$ws.Range("D15").Value = '98.225.91'
$ws.Range("E15").Value = '  -0.46%  '

$ws.Range("D16").Value = '3.983.28'
$ws.Range("E16").Value = '  -0.13%  '

$ws.Range("D17").Value = '5.54'
$ws.Range("E17").Value = '  -0.67%  '

$ws.Range("D18").Value = '3.336.56'
$ws.Range("E18").Value = '  -1.60%  '

$ws.Range("E19").Value = '  +18.22%  '

$ws.Range("D20").Value = '16.68'
$ws.Range("E20").Value = '  +9.07%  '

$ws.Range("D21").Value = '531.17'
$ws.Range("E21").Value = '  +7.28%  '

$ws.Range("D22").Value = '3.55'
$ws.Range("E22").Value = '  -2.68%  '

$ws.Range("D23").Value = '10.06'
$ws.Range("E23").Value = '  +6.92%  '

$ws.Range("D24").Value = '0.0000210'
$ws.Range("E24").Value = '  -1.19%  '

$ws.Range("D25").Value = '0.418'
$ws.Range("E25").Value = '  +46.33%  '

$ws.Range("D26").Value = '102.11'
$ws.Range("E26").Value = '  +14.52%  '

$ws.Range("D27").Value = '6.07'
$ws.Range("E27").Value = '  +5.37%  '

$ws.Range("D28").Value = '12.75'
$ws.Range("E28").Value = '  +5.31%  '

$ws.Range("D29").Value = '3.518.89'
$ws.Range("E29").Value = '  -0.01%  '

$ws.Range("E30").Value = '  +12.05%  '

$ws.Range("D31").Value = '0.998'
$ws.Range("E31").Value = '  +0.02%  '

$ws.Range("D32").Value = '10.85'
$ws.Range("E32").Value = '  +12.36%  '

$ws.Range("D33").Value = '0.188'
$ws.Range("E33").Value = '  -4.87%  '

$ws.Range("E34").Value = '  +0.20%  '

$ws.Range("D35").Value = '28.95'
$ws.Range("E35").Value = '  +3.26%  '

$ws.Range("D36").Value = '0.512'
$ws.Range("E36").Value = '  +8.86%  '

$ws.Range("D37").Value = '7.78'
$ws.Range("E37").Value = '  +5.02%  '

$ws.Range("E38").Value = '  +2.78%  '

$ws.Range("D39").Value = '2.05'
$ws.Range("E39").Value = '  +2.50%  '

$ws.Range("D40").Value = '519.01'
$ws.Range("E40").Value = '  +2.27%  '

$ws.Range("E41").Value = '  -0.67%  '

$ws.Range("B42").Value = 'Fetch.AI'
$ws.Range("C42").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D42").Value = '1.31'
$ws.Range("E42").Value = '  +2.63%  '

$ws.Range("B43").Value = 'MantraDAO'
$ws.Range("C43").Value = 'https://coinranking.com/coin/cTdD8lD-6+mantradao-om'
$ws.Range("D43").Value = '3.85'
$ws.Range("E43").Value = '  +2.78%  '

$ws.Range("D44").Value = '0.806'
$ws.Range("E44").Value = '  +2.35%  '

$ws.Range("D45").Value = '3.30'
$ws.Range("E45").Value = '  -3.09%  '

$ws.Range("D47").Value = '0.0390'
$ws.Range("E47").Value = '  +19.41%  '

$ws.Range("D48").Value = '163.94'
$ws.Range("E48").Value = '  +2.14%  '

$ws.Range("D49").Value = '2.01'
$ws.Range("E49").Value = '  +3.07%  '

$ws.Range("D50").Value = '7.69'
$ws.Range("E50").Value = '  +16.47%  '

$ws.Range("D51").Value = '49.44'
$ws.Range("E51").Value = '  +5.89%  '
